$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.754.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.41%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.917.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.50%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'241.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.74%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.14%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4923"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.31%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3004"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.74%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06782"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.919.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.81%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'17.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.37%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'88.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6774"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.04%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'30.737.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.42%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000007995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.40%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'13.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.54%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.14%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.165.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.458"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +13.19%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.02%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'199.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +6.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.376"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.84%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.700"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.82%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'162.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +4.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.85%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.968"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.23%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.474"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.46%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.366"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.75%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09174"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.082"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05328"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.57%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7470"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.37%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.01862"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.08%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.730"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.9325"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.18%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.097"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.34%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.4522"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Quant"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'107.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.37%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.970"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Aave"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'72.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +24.83%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'PaxDollar"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.003"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Algorand"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.1402"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.72%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'7.743"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.91%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'9.169"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.51%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Elrond"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'35.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +6.49%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05910"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.00%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4060"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.25%  "
$ws.Range("E51").Style = "Normal"
